# Update scripts with new TPM data (Cxcl13-Cxcr5, YoungD0, FAPs-only rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rerun with updated TPM values only produced results for the FAPs
# sending cluster (rows 2-4); the former MuSCs sending-cluster rows
# (5-7) are no longer present, so remove them.
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2 (FAPs -> Cxcl13/Cxcr5 -> ECs) numeric refresh
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1786973333333333
$ws.Range("N2").Value = 0.536092
$ws.Range("O2").Value = 0.1393273670217259
$ws.Range("P2").Value = 0.1393273670217259
$ws.Range("Q2").Value = 0.7165490851062223
$ws.Range("R2").Value = 6.448941765956
$ws.Range("S2").Value = 0.1393273670217259
$ws.Range("T2").Value = 0.1393273670217259

# Row 3 (FAPs -> Cxcl13/Cxcr5 -> FAPs) numeric refresh
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.5248434460452502
$ws.Range("P3").Value = 0.5248434460452502
$ws.Range("S3").Value = 0.5248434460452502
$ws.Range("T3").Value = 0.5248434460452502

# Row 4 (FAPs -> Cxcl13/Cxcr5 -> MuSCs) numeric refresh
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.430725
$ws.Range("N4").Value = 1.292175
$ws.Range("O4").Value = 0.3358291869330239
$ws.Range("P4").Value = 0.3358291869330239
$ws.Range("Q4").Value = 1.727141636225
$ws.Range("R4").Value = 15.544274726025
$ws.Range("S4").Value = 0.3358291869330239
$ws.Range("T4").Value = 0.3358291869330239
